$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '70.438.39'
$ws.Range("E2").Value = '  +2.23%  '
$ws.Range("D3").Value = '3.809.25'
$ws.Range("E3").Value = '  +0.96%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '''681.35'
$ws.Range("E5").Value = '  +8.40%  '
$ws.Range("D6").Value = '''171.29'
$ws.Range("E6").Value = '  +3.82%  '
$ws.Range("D7").Value = '3.807.84'
$ws.Range("E7").Value = '  +1.03%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("E9").Value = '  +0.79%  '
$ws.Range("E10").Value = '  +1.66%  '
$ws.Range("E11").Value = '  +7.00%  '
$ws.Range("E12").Value = '  +0.80%  '
$ws.Range("E13").Value = '  +0.61%  '
$ws.Range("D14").Value = '''35.98'
$ws.Range("E14").Value = '  +2.40%  '
$ws.Range("D15").Value = '4.454.74'
$ws.Range("E15").Value = '  +1.04%  '
$ws.Range("D16").Value = '3.811.22'
$ws.Range("E16").Value = '  +1.09%  '
$ws.Range("D17").Value = '70.538.79'
$ws.Range("E18").Value = '  +0.78%  '
$ws.Range("E19").Value = '  +2.24%  '
$ws.Range("D20").Value = '''0.114'
$ws.Range("E20").Value = '  +0.72%  '
$ws.Range("D21").Value = '''11.26'
$ws.Range("E21").Value = '  +18.21%  '
$ws.Range("D22").Value = '''477.23'
$ws.Range("E22").Value = '  +2.46%  '
$ws.Range("D23").Value = '''0.714'
$ws.Range("E23").Value = '  +1.20%  '
$ws.Range("D24").Value = '''83.41'
$ws.Range("E24").Value = '  +0.56%  '
$ws.Range("E25").Value = '  -1.55%  '
$ws.Range("D26").Value = '''12.25'
$ws.Range("E26").Value = '  +2.22%  '
$ws.Range("E27").Value = '  +3.24%  '
$ws.Range("E28").Value = '  -1.10%  '
$ws.Range("E29").Value = '  +0.05%  '
$ws.Range("D30").Value = '3.960.96'
$ws.Range("E30").Value = '  +0.96%  '
$ws.Range("D31").Value = '''2.93'
$ws.Range("E31").Value = '  +9.76%  '
$ws.Range("E32").Value = '  +3.00%  '
$ws.Range("D33").Value = '''7.40'
$ws.Range("E33").Value = '  +4.08%  '
$ws.Range("D34").Value = '''29.60'
$ws.Range("E35").Value = '  +5.36%  '
$ws.Range("B36").Value = 'Aptos'
$ws.Range("C36").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D36").Value = '''9.14'
$ws.Range("E36").Value = '  +2.38%  '
$ws.Range("B37").Value = 'Binance-PegBSC-USD'
$ws.Range("C37").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D37").Value = '''0.998'
$ws.Range("E37").Value = '  -0.18%  '
$ws.Range("D38").Value = '3.759.82'
$ws.Range("E38").Value = '  +0.95%  '
$ws.Range("E39").Value = '  +1.54%  '
$ws.Range("D40").Value = '''3.39'
$ws.Range("E40").Value = '  +3.12%  '
$ws.Range("D41").Value = '''5.95'
$ws.Range("E41").Value = '  +2.51%  '
$ws.Range("E42").Value = '  -0.25%  '
$ws.Range("E43").Value = '  +0.06%  '
$ws.Range("D44").Value = '''2.14'
$ws.Range("E44").Value = '  +13.33%  '
$ws.Range("D46").Value = '''46.23'
$ws.Range("E46").Value = '  +7.52%  '
$ws.Range("D47").Value = '''160.02'
$ws.Range("E47").Value = '  +2.89%  '
$ws.Range("D48").Value = '''1.46'
$ws.Range("E48").Value = '  +7.67%  '
$ws.Range("D49").Value = '''48.15'
$ws.Range("E49").Value = '  +3.29%  '
$ws.Range("D50").Value = '''0.000295'
$ws.Range("E50").Value = '  +8.74%  '
$ws.Range("E51").Value = '  +1.82%  '
